# Insert a new data row at row 181 (pushing existing rows 181..272 down to 182..273)
# and populate it with the new "Granada" price observation for Vega Modelo de Temuco.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A181").EntireRow.Insert()

$ws.Range("A181").Value = 10
$ws.Range("B181").Value = "Vega Modelo de Temuco"
$ws.Range("C181").Value = "La Araucanía"
$ws.Range("D181").Value = 45134
$ws.Range("E181").Value = 9
$ws.Range("F181").Value = "Fruta"
$ws.Range("G181").Value = 100104
$ws.Range("H181").Value = "Frutos de pepita"
$ws.Range("I181").Value = 100104001
$ws.Range("J181").Value = "Granada"
$ws.Range("K181").Value = "Wonderfull"
$ws.Range("L181").Value = "Primera"
$ws.Range("M181").Value = 250
$ws.Range("N181").Value = 16000
$ws.Range("O181").Value = 16000
$ws.Range("P181").Value = 16000
$ws.Range("Q181").Value = "$/bandeja 10 kilos"
$ws.Range("R181").Value = "Provincia de Limarí"
$ws.Range("S181").Value = 1600
$ws.Range("T181").Value = 10
